$d = $word.ActiveDocument

# Locate the paragraph whose full text is exactly "Required Skill Sets:"
# (there is another, unrelated, occurrence of "Required Skill Sets" in the
# Heading2 title a couple of paragraphs above, so match the colon too).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n", [char]7)
    if ($t -eq "Required Skill Sets:") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Range covering just the run's text, excluding the paragraph mark.
    $full = $target.Range
    $r = $d.Range($full.Start, $full.End - 1)

    # Swap the wording but keep the bold formatting.
    $r.Text = "Skills you will gain:"
    $r.Bold = 1

    # Force the trailing colon into its own run (matching the target
    # markup, which has "Skills you will gain" and ":" as separate runs)
    # by toggling its bold state off and back on, which splits the run
    # without altering the visible formatting.
    $colon = $d.Range($r.End - 1, $r.End)
    $colon.Bold = 0
    $colon.Bold = 1
}
